# otros.xlsx -- add the funerary-goods rows to the "instrumentos" sheet
# (Tumba con monolito / Mausoleo / Ataud / Cripta / Monticulo funerario /
#  Sarcofago / Urna funeraria), renumber the trailing blank spacer rows,
# and restore the sheet view/selection left by the author's last save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft = -4131

# ---------------------------------------------------------------------
# 1. New data rows 16-22 (nombre, peso, precio, descripcion)
# ---------------------------------------------------------------------

# Rows 20 and 21 already existed in the sheet as blank, left-aligned
# spacer cells (column D only, no value). Row 20 keeps that same
# left-aligned style once it gets real data below, but row 21 needs to
# go back to the plain/default style -- so its stale formatting is
# stripped first instead of leaking into the new content.
$ws.Cells.Item(21, 4).Clear()

$newRows = @(
    @{ Row = 16; Nombre = "Tumba con monolito";  Peso = "-"; Precio = "50L";     Desc = "Con capacidad para uno o dos difuntos."; Align = $false },
    @{ Row = 17; Nombre = "Mausoleo";             Peso = "-"; Precio = "25L";     Desc = "Para un solo difunto.";                  Align = $false },
    @{ Row = 18; Nombre = "Ataud";                Peso = "-"; Precio = "2L";      Desc = "De madera y sin adornos.";               Align = $false },
    @{ Row = 19; Nombre = "Cripta";               Peso = "-"; Precio = "100L";    Desc = "Cámara funeraria subterranea.";          Align = $true  },
    @{ Row = 20; Nombre = "Montículo funerario";  Peso = "-"; Precio = "20-100L"; Desc = "Montículo con cripta subterranea.";      Align = $true  },
    @{ Row = 21; Nombre = "Sarcófago";             Peso = "-"; Precio = "20-100";  Desc = "De madera o piedra tallada, con o sin incrustaciones y policromía."; Align = $false },
    @{ Row = 22; Nombre = "Urna funeraria";       Peso = "-"; Precio = "5C-10L";  Desc = "Según material y calidad.";             Align = $true  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Nombre
    $ws.Cells.Item($row, 2).Value = $r.Peso
    $ws.Cells.Item($row, 3).Value = $r.Precio
    $ws.Cells.Item($row, 4).Value = $r.Desc
    if ($r.Align) {
        $ws.Cells.Item($row, 4).HorizontalAlignment = $xlLeft
    }
}

# ---------------------------------------------------------------------
# 2. Drop the blank spacer rows that no longer apply (their row numbers
#    shift now that rows 16-22 hold real data) and recreate the blank
#    placeholder rows further down at their new row numbers.
# ---------------------------------------------------------------------
$ws.Cells.Item(25, 4).Clear()
$ws.Cells.Item(30, 4).Clear()
$ws.Cells.Item(39, 4).Clear()
$ws.Cells.Item(45, 4).Clear()

$ws.Cells.Item(29, 4).HorizontalAlignment = $xlLeft
$ws.Cells.Item(36, 4).HorizontalAlignment = $xlLeft
$ws.Cells.Item(44, 4).HorizontalAlignment = $xlLeft

# ---------------------------------------------------------------------
# 3. Restore the view (zoom + active selection) from the author's save.
# ---------------------------------------------------------------------
$ws.Range("A23").Select()
$excel.ActiveWindow.Zoom = 115
